# Rename the two repeated logo pictures (Pearson logo in the footers,
# BTec logo in the headers) so each inline picture's shape Name matches
# the image it actually displays instead of the other logo's filename.
#   Pearson logo (descr ends in PearsonLogo.png) -> "image2.png"
#   BTec logo   (descr "BTec_Logo-Orange")       -> "image1.jpg"
#
# This only touches the cosmetic <wp:docPr name="..."> identifier (what
# Word calls the picture's "Name" in the Selection Pane) - the
# embedded picture data / relationships are untouched.

$d = $word.ActiveDocument

function Rename-LogoShapes($range) {
    if ($null -eq $range) { return }
    $shapes = $range.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $descr = $shp.AlternativeText
        if ($descr -like "*PearsonLogo.png") {
            $shp.Name = "image2.png"
        } elseif ($descr -eq "BTec_Logo-Orange") {
            $shp.Name = "image1.jpg"
        }
    }
}

foreach ($sec in $d.Sections) {
    for ($hIdx = 1; $hIdx -le 3; $hIdx++) {
        $hdr = $sec.Headers($hIdx)
        if ($hdr.Exists) {
            Rename-LogoShapes $hdr.Range
        }
    }
    for ($fIdx = 1; $fIdx -le 3; $fIdx++) {
        $ftr = $sec.Footers($fIdx)
        if ($ftr.Exists) {
            Rename-LogoShapes $ftr.Range
        }
    }
}

Write-Output "Renamed logo inline shapes."
